$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Steps")

# Extend the formatted block (rows 11-18) by copying the format of the last
# existing data row (row 10) down into the new rows so the new cells share
# the same style (fill/border) as the rest of the table.
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update existing rows 7-10 ---
$ws.Range("C7").Value = "waitForElementPresent"

$ws.Range("C8").Value = "click"
$ws.Range("E8").Value = ""

$ws.Range("C9").Value = "type"
$ws.Range("D9").Value = "search_box"

$ws.Range("C10").Value = "click"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = "apple"

# --- Row 11 (reuses existing shared strings only) ---
$ws.Range("A11").Value = "searchListingPageTestCases"
$ws.Range("C11").Value = "addProduct"
$ws.Range("D11").Value = "search_list_projectNames|search_list_addBtn"
$ws.Range("E11").Value = "Apple - Royal Gala"

# --- Row 12 (introduces verifyText, productName_text) ---
$ws.Range("A12").Value = "searchListingPageTestCases"
$ws.Range("C12").Value = "verifyText"
$ws.Range("D12").Value = "productName_text"
$ws.Range("E12").Value = "Apple - Royal Gala"

# --- Row 13 (introduces product_add_btn) ---
$ws.Range("A13").Value = "searchListingPageTestCases"
$ws.Range("C13").Value = "click"
$ws.Range("D13").Value = "product_add_btn"

# Introduce "search_list_footerpage" before "goBack"/"Green Apple" so the
# shared-string table ends up in the same order as the target workbook.
$ws.Range("D18").Value = "search_list_footerpage"

# --- Row 14 (introduces goBack) ---
$ws.Range("A14").Value = "searchListingPageTestCases"
$ws.Range("C14").Value = "goBack"

# --- Row 15 (introduces Green Apple) ---
$ws.Range("A15").Value = "searchListingPageTestCases"
$ws.Range("C15").Value = "addProduct"
$ws.Range("D15").Value = "search_list_projectNames|search_list_addBtn"
$ws.Range("E15").Value = "Green Apple"

# --- Row 16 ---
$ws.Range("A16").Value = "searchListingPageTestCases"
$ws.Range("C16").Value = "verifyText"
$ws.Range("D16").Value = "productName_text"
$ws.Range("E16").Value = "Green Apple"

# --- Row 17 ---
$ws.Range("A17").Value = "searchListingPageTestCases"
$ws.Range("C17").Value = "click"
$ws.Range("D17").Value = "product_add_btn"

# --- Row 18 (remaining cells; D18 already set above) ---
$ws.Range("A18").Value = "searchListingPageTestCases"
$ws.Range("C18").Value = "click"

# Update the sheet view to match the author's final cursor/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("E12").Select()
